# Apply updated cryptocurrency price/volume data to Sheet1 (columns D and E).
# Matches the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.827.62"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "1.736.71"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.81"
$ws.Range("E5").Value = "  +2.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5113"
$ws.Range("E7").Value = "  -1.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2729"
$ws.Range("E8").Value = "  -1.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06107"
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("D10").Value = "1.738.01"
$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6370"
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.578"
$ws.Range("E14").Value = "  +1.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.06"
$ws.Range("E15").Value = "  +0.47%  "

$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").Value = "25.839.58"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.69"
$ws.Range("E19").Value = "  +2.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006738"
$ws.Range("E20").Value = "  +1.67%  "

$ws.Range("D21").Value = "1.961.22"
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.249"
$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.662"
$ws.Range("E23").Value = "  -0.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.213"
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.75"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.09"
$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("E28").Value = "  -1.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.49"
$ws.Range("E29").Value = "  +3.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.980"
$ws.Range("E30").Value = "  +8.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08305"
$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.642"
$ws.Range("E32").Value = "  +4.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04562"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.664"
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9843"
$ws.Range("E35").Value = "  +0.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6153"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.682"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01604"
$ws.Range("E38").Value = "  +2.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.922"
$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9998"
$ws.Range("E40").Value = "  +0.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.99"
$ws.Range("E41").Value = "  -1.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3827"
$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7349"
$ws.Range("E43").Value = "  +2.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.945"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1119"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05255"
$ws.Range("E46").Value = "  -1.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.135"
$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.73"
$ws.Range("E48").Value = "  +3.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.38"
$ws.Range("E49").Value = "  +1.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.569"
$ws.Range("E50").Value = "  -1.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3400"
$ws.Range("E51").Value = "  +0.81%  "
